# Apply the "new item uploaded" edit to the DaySale workbook:
#  1. Insert a new product row ("PROPAMETHONE TOP. CREAM. 20 GM") into the
#     shortage table, keeping it in alphabetical order right before
#     "PROSTRIDE 5MG 30 CAPS." (i.e. as the new row 44, pushing the rest of
#     the table down by one row).
#  2. Update the printed timestamp string (7:07 PM -> 7:10 PM).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new row ------------------------------------------------
# Row 44 currently holds "PROSTRIDE 5MG 30 CAPS." (item #38). Shift it (and
# everything below it) down by one row so a new row 44 is free.
$ws.Rows(44).Insert(-4121, 0)

# Paint the new row with the same formatting used by the rest of the table
# (border / fill / font / number formats) by copying it from the row that
# used to be "row 44" and is now "row 45".
$ws.Range("A45:Q45").Copy()
$ws.Range("A44:Q44").PasteSpecial(-4122)
$ws.Rows(44).RowHeight = $ws.Rows(45).RowHeight()

# Recreate the merged cells for the new row (matching every other data row
# in the table: A:B, C:G, H:K, L:M, N:O).
$ws.Range("A44:B44").Merge()
$ws.Range("C44:G44").Merge()
$ws.Range("H44:K44").Merge()
$ws.Range("L44:M44").Merge()
$ws.Range("N44:O44").Merge()

# Fill in the new item's data.
$ws.Range("A44").Value = 38
$ws.Range("C44").Value = "PROPAMETHONE TOP. CREAM. 20 GM"
$ws.Range("H44").Value = "0:0"
$ws.Range("L44").Value = "1"
$ws.Range("N44").Value = "30.00"
$ws.Range("P44").Value = "30.0000"
$ws.Range("Q44").Value = "1:0"

# --- 2. Bump the printed "generated at" timestamp --------------------------
# This label lives in the footer row, which used to be row 71 and (after the
# row insert above) is now row 72. Look it up by value instead of a hard
# coded row number so the script stays correct even if that ever changes.
$found = $false
for ($r = 1; $r -le $ws.UsedRange.Rows.Count(); $r++) {
    $v = $ws.Cells.Item($r, 1).Value()
    if ($v -eq "Monday, 9 June, 2025 7:07 PM") {
        $ws.Cells.Item($r, 1).Value = "Monday, 9 June, 2025 7:10 PM"
        $found = $true
    }
}
if (-not $found) {
    $ws.Range("A72").Value = "Monday, 9 June, 2025 7:10 PM"
}
